$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 199, shifting existing rows 199-231 down to 200-232.
$ws.Rows(199).Insert()

# Populate the newly inserted row 199 with the new data entry.
$ws.Range("A199").Value = 3
$ws.Range("B199").Value = "Femacal de La Calera"
$ws.Range("C199").Value = "Coquimbo"
$ws.Range("D199").Value = 44946
$ws.Range("E199").Value = 5
$ws.Range("F199").Value = 100112030
$ws.Range("G199").Value = "Poroto granado"
$ws.Range("H199").Value = "Sin especificar"
$ws.Range("I199").Value = "Primera"
$ws.Range("J199").Value = 38
$ws.Range("K199").Value = 40000
$ws.Range("L199").Value = 40000
$ws.Range("M199").Value = 40000
$ws.Range("N199").Value = "`$/saco 25 kilos"
$ws.Range("O199").Value = "Provincia de Limarí"
$ws.Range("P199").Value = 1600
$ws.Range("Q199").Value = 25
$ws.Range("R199").Value = "Hortaliza"
